$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D16").Value = "Grad-CAM: Why did you say that? 내용 정리 [XAI-14]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/163"

$ws.Range("D41").Value = "클라우드 서비스와 오픈소스 라이선싱"
$ws.Range("E41").Value = "http://cloudinsight.net/cloud/%ed%81%b4%eb%9d%bc%ec%9a%b0%eb%93%9c-%ec%84%9c%eb%b9%84%ec%8a%a4%ec%99%80-%ec%98%a4%ed%94%84%ec%86%8c%ec%8a%a4-%eb%9d%bc%ec%9d%b4%ec%84%a0%ec%8b%b1/"
